# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 26 de Mayo de 2020 a las 06:35"

# "Haiti" moves up the ranking to just after "Hong Kong" (row 108), pushing
# "Mali", "Tunez", "Letonia", "Guinea Ecuatorial" and "Albania" down one row
# each (rows 109-113). Refreshed case numbers are applied at the same time.

# Row 100 (Kirguistan) - refreshed stats only, ranking unchanged
$ws.Range("B100").Value = 1468
$ws.Range("C100").Value = 35
$ws.Range("D100").Value = 1015
$ws.Range("E100").Value = 437

# Row 108: Mali -> Haiti, with Haiti's fresh stats
$ws.Range("A108").Value = "Haiti"
$ws.Range("B108").Value = 1063
$ws.Range("C108").Value = 105
$ws.Range("D108").Value = 22
$ws.Range("E108").Value = 1010
$ws.Range("G108").Value = 4
$ws.Range("H108").Value = 31

# Row 109: Tunez -> Mali (keeps Mali's previous stats)
$ws.Range("A109").Value = "Mali"
$ws.Range("B109").Value = 1059
$ws.Range("C109").Value = 0
$ws.Range("D109").Value = 604
$ws.Range("E109").Value = 388
$ws.Range("H109").Value = 67

# Row 110: Letonia -> Tunez (keeps Tunez's previous stats)
$ws.Range("A110").Value = "Tunez"
$ws.Range("B110").Value = 1051
$ws.Range("C110").Value = 0
$ws.Range("D110").Value = 919
$ws.Range("E110").Value = 84
$ws.Range("H110").Value = 48

# Row 111: Guinea Ecuatorial -> Letonia (keeps Letonia's previous stats)
$ws.Range("A111").Value = "Letonia"
$ws.Range("B111").Value = 1049
$ws.Range("C111").Value = 0
$ws.Range("D111").Value = 712
$ws.Range("E111").Value = 315
$ws.Range("H111").Value = 22

# Row 112: Albania -> Guinea Ecuatorial (keeps Guinea Ecuatorial's previous stats)
$ws.Range("A112").Value = "Guinea Ecuatorial"
$ws.Range("B112").Value = 1043
$ws.Range("C112").Value = 0
$ws.Range("D112").Value = 165
$ws.Range("E112").Value = 866
$ws.Range("H112").Value = 12

# Row 113: Haiti -> Albania (keeps Albania's previous stats)
$ws.Range("A113").Value = "Albania"
$ws.Range("B113").Value = 1004
$ws.Range("C113").Value = 0
$ws.Range("D113").Value = 795
$ws.Range("E113").Value = 177
$ws.Range("H113").Value = 32
